$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.790418148040771
$ws.Range("B1").Value = 4.27653980255127
$ws.Range("C1").Value = 1.646914839744568
$ws.Range("D1").Value = 0.8684206604957581
$ws.Range("E1").Value = 0.4686254262924194
